$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.774.76"
$ws.Cells.Item(2, 5).Value = "  -0.11%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.634.55"
$ws.Cells.Item(3, 5).Value = "  -0.04%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.28%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'215.78"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.30%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.82%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.27%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.32%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -1.26%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'19.60"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.48%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0793"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.73%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.02%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.860.22"
$ws.Cells.Item(13, 5).Value = "  -0.09%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.638.50"
$ws.Cells.Item(14, 5).Value = "  -0.03%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.561"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.54%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.0₃0762"
$ws.Cells.Item(16, 5).Value = "  -0.74%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.15%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.800.72"
$ws.Cells.Item(18, 5).Value = "  -0.10%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.19%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'4.47"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.59%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'192.48"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.69%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'9.97"
$ws.Cells.Item(22, 4).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 5).Value = "  +2.40%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +3.68%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.17%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'141.82"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.19%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.96%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'6.90"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.94%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'15.49"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.10%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.23%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.38%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.33"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.04%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.83%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.93%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.17%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.19%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "1.131.66"
$ws.Cells.Item(37, 5).Value = "  +1.58%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "'0.545"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.47%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).Value = "'2.51"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.25%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -1.42%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.02%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "mCoin"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Cells.Item(42, 4).Value = "'2.53"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.49%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).Value = "'5.56"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.31%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Quant"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(44, 4).Value = "'100.72"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.09%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(45, 4).Value = "'0.799"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.39%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "1.769.37"
$ws.Cells.Item(46, 5).Value = "  -0.21%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(47, 4).Value = "0.0₆0112"
$ws.Cells.Item(47, 5).Value = "  +0.31%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).Value = "'55.39"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.24%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).Value = "'0.416"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.03%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.0503"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.08%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).Value = "'1.42"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +4.02%  "
